$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For numeric-looking text values (e.g. "0.546"), set the cell
# NumberFormat to Text first so Excel keeps the value as a string
# instead of auto-converting it to a floating point number.

$ws.Range("D2").Value = "34.106.18"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "1.785.83"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("E4").Value = "  +0.35%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.25"
$ws.Range("E5").Value = "  -0.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.546"
$ws.Range("E6").Value = "  -0.91%  "
$ws.Range("E7").Value = "  +0.38%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.91"
$ws.Range("E8").Value = "  -3.30%  "
$ws.Range("E9").Value = "  +0.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0689"
$ws.Range("E10").Value = "  -3.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0944"
$ws.Range("E11").Value = "  +1.00%  "
$ws.Range("D12").Value = "2.042.00"
$ws.Range("E12").Value = "  -0.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.22"
$ws.Range("E13").Value = "  +0.64%  "
$ws.Range("D14").Value = "1.796.06"
$ws.Range("E14").Value = "  +0.50%  "
$ws.Range("D15").Value = "34.037.21"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.618"
$ws.Range("E16").Value = "  -1.32%  "
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.93"
$ws.Range("E18").Value = "  -0.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.49"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").Value = "0.0₃0777"
$ws.Range("E20").Value = "  -1.54%  "
$ws.Range("E21").Value = "  +0.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.81"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.09"
$ws.Range("E23").Value = "  -0.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.04"
$ws.Range("E24").Value = "  -2.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.53"
$ws.Range("E25").Value = "  +0.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.14"
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.30"
$ws.Range("E28").Value = "  +0.42%  "
$ws.Range("E29").Value = "  +0.52%  "
$ws.Range("E30").Value = "  -0.79%  "
$ws.Range("E31").Value = "  +0.33%  "
$ws.Range("E32").Value = "  -0.71%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.61"
$ws.Range("E33").Value = "  +2.21%  "
$ws.Range("E34").Value = "  -0.97%  "
$ws.Range("D35").Value = "1.453.44"
$ws.Range("E35").Value = "  +3.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.645"
$ws.Range("E36").Value = "  -1.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0192"
$ws.Range("E37").Value = "  +1.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.40"
$ws.Range("E38").Value = "  +8.06%  "
$ws.Range("E39").Value = "  -0.98%  "
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "80.06"
$ws.Range("E40").Value = "  +1.27%  "
$ws.Range("B41").Value = "HuobiToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.37"
$ws.Range("E41").Value = "  +0.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.918"
$ws.Range("E42").Value = "  -0.35%  "
$ws.Range("E43").Value = "  -0.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.44"
$ws.Range("E44").Value = "  +2.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0509"
$ws.Range("E45").Value = "  +2.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.03"
$ws.Range("E46").Value = "  +3.43%  "
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("D48").Value = "0.0₆0138"
$ws.Range("E48").Value = "  -0.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "107.46"
$ws.Range("E49").Value = "  -1.68%  "
$ws.Range("D50").Value = "1.943.75"
$ws.Range("E50").Value = "  -0.14%  "
$ws.Range("E51").Value = "  +0.42%  "
